$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 266 ("Valencia"/44188 record),
# shifting the existing rows 266-280 down to 268-282.
$ws.Rows.Item(266).Insert()
$ws.Rows.Item(266).Insert()

# Populate the two newly inserted rows with the new "New Hall" records.
$ws.Cells.Item(266, 1).Value = 11
$ws.Cells.Item(266, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(266, 3).Value = "Bíobío"
$ws.Cells.Item(266, 4).Value = 44714
$ws.Cells.Item(266, 5).Value = 8
$ws.Cells.Item(266, 6).Value = "Fruta"
$ws.Cells.Item(266, 7).Value = 100102
$ws.Cells.Item(266, 8).Value = "Cítricos"
$ws.Cells.Item(266, 9).Value = 100102005
$ws.Cells.Item(266, 10).Value = "Naranja"
$ws.Cells.Item(266, 11).Value = "New Hall"
$ws.Cells.Item(266, 12).Value = "Primera"
$ws.Cells.Item(266, 13).Value = 100
$ws.Cells.Item(266, 14).Value = 8000
$ws.Cells.Item(266, 15).Value = 9000
$ws.Cells.Item(266, 16).Value = 8500
$ws.Cells.Item(266, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(266, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(266, 19).Value = 567
$ws.Cells.Item(266, 20).Value = 15

$ws.Cells.Item(267, 1).Value = 11
$ws.Cells.Item(267, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(267, 3).Value = "Bíobío"
$ws.Cells.Item(267, 4).Value = 44714
$ws.Cells.Item(267, 5).Value = 8
$ws.Cells.Item(267, 6).Value = "Fruta"
$ws.Cells.Item(267, 7).Value = 100102
$ws.Cells.Item(267, 8).Value = "Cítricos"
$ws.Cells.Item(267, 9).Value = 100102005
$ws.Cells.Item(267, 10).Value = "Naranja"
$ws.Cells.Item(267, 11).Value = "New Hall"
$ws.Cells.Item(267, 12).Value = "Segunda"
$ws.Cells.Item(267, 13).Value = 50
$ws.Cells.Item(267, 14).Value = 7000
$ws.Cells.Item(267, 15).Value = 7000
$ws.Cells.Item(267, 16).Value = 7000
$ws.Cells.Item(267, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(267, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(267, 19).Value = 467
$ws.Cells.Item(267, 20).Value = 15
